$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$headers = @(
  "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310",
  "Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310","Bedingung_FV2310","diff",
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404",
  "Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404","Bedingung_FV2404"
)
for ($c = 1; $c -le 21; $c++) {
  $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U84")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (split + frozen pane) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
